$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.341.28"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").Value = "3.485.39"
$ws.Range("E3").Value = "  -1.27%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.98%  "

$ws.Range("D7").Value = "3.482.90"
$ws.Range("E7").Value = "  -1.30%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +0.67%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.144"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.59"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.77%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.429"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.70%  "

$ws.Range("E13").Value = "  -2.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.48%  "

$ws.Range("D15").Value = "4.068.26"
$ws.Range("E15").Value = "  -1.28%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.495.02"
$ws.Range("E16").Value = "  -0.86%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "67.231.45"
$ws.Range("E17").Value = "  -0.41%  "

$ws.Range("E18").Value = "  -0.83%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.77%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "444.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.626"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.99%  "

$ws.Range("D26").Value = "3.621.59"
$ws.Range("E26").Value = "  -1.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000126"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.44%  "

$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.24%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.84%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.62%  "

$ws.Range("E32").Value = "  +3.20%  "

$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.46%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.03%  "

$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").Value = "3.477.28"
$ws.Range("E37").Value = "  -1.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.98%  "

$ws.Range("E40").Value = "  +6.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "177.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0894"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.33%  "

$ws.Range("E45").Value = "  +0.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.53"
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.33%  "

$ws.Range("E51").Value = "  -0.05%  "
